# Scheduled runner update: refresh market-derived profit figures (Anima_Profits)
# across all Sheets, row by row, preserving original cell layout.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 10052.071
$ws.Range("J112").Value = 10052.071
$ws.Range("L112").Value = 30156.213
$ws.Range("N112").Value = -32372.213
$ws.Range("H113").Value = 2417.8262
$ws.Range("I113").Value = 1959.1666
$ws.Range("K113").Value = 1959.1666
$ws.Range("M113").Value = 1294.8334
$ws.Range("H116").Value = 6385.8335
$ws.Range("I116").Value = 12370.5
$ws.Range("J116").Value = 2111.0715
$ws.Range("K116").Value = 12370.5
$ws.Range("L116").Value = 2111.0715
$ws.Range("M116").Value = -8928.5
$ws.Range("N116").Value = -8995.0715
$ws.Range("H132").Value = 2781.9824
$ws.Range("I132").Value = 2634.7844
$ws.Range("K132").Value = 7904.3532
$ws.Range("M132").Value = -5374.3532
$ws.Range("H137").Value = 1101.3296
$ws.Range("I137").Value = 1001.39215
$ws.Range("J137").Value = 1239.081
$ws.Range("K137").Value = 3004.17645
$ws.Range("L137").Value = 3717.242999999999
$ws.Range("M137").Value = -454.1764499999999
$ws.Range("N137").Value = -8817.242999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 651148.3
$ws.Range("I32").Value = 740667.1
$ws.Range("J32").Value = 17056.416
$ws.Range("K32").Value = 740667.1
$ws.Range("L32").Value = 17056.416
$ws.Range("M32").Value = -740380.1
$ws.Range("N32").Value = -17630.416
$ws.Range("H61").Value = 7577679
$ws.Range("I61").Value = 8773602
$ws.Range("J61").Value = 3498
$ws.Range("K61").Value = 8773602
$ws.Range("L61").Value = 3498
$ws.Range("M61").Value = -8773390
$ws.Range("N61").Value = -3922
$ws.Range("H74").Value = 1445
$ws.Range("I74").Value = 761.3684
$ws.Range("K74").Value = 761.3684
$ws.Range("M74").Value = 112.6316
$ws.Range("H77").Value = 1445
$ws.Range("I77").Value = 761.3684
$ws.Range("K77").Value = 3806.842
$ws.Range("M77").Value = 561.1580000000004
$ws.Range("H132").Value = 2870
$ws.Range("I132").Value = 3048.2068
$ws.Range("K132").Value = 9144.6204
$ws.Range("M132").Value = -6614.6204
$ws.Range("H136").Value = 7577679
$ws.Range("I136").Value = 8773602
$ws.Range("J136").Value = 3498
$ws.Range("K136").Value = 26320806
$ws.Range("L136").Value = 10494
$ws.Range("M136").Value = -26318256
$ws.Range("N136").Value = -15594

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 246.21213
$ws.Range("I22").Value = 246.21213
$ws.Range("K22").Value = 246.21213
$ws.Range("M22").Value = -73.21213
$ws.Range("H107").Value = 112522.336
$ws.Range("I107").Value = 201060.2
$ws.Range("J107").Value = 1850
$ws.Range("K107").Value = 201060.2
$ws.Range("L107").Value = 1850
$ws.Range("M107").Value = -199140.2
$ws.Range("N107").Value = -5690
$ws.Range("H134").Value = 2327.0588
$ws.Range("I134").Value = 2273.7932
$ws.Range("K134").Value = 6821.3796
$ws.Range("M134").Value = -4286.3796

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 620
$ws.Range("I22").Value = 565
$ws.Range("J22").Value = 675
$ws.Range("K22").Value = 565
$ws.Range("L22").Value = 675
$ws.Range("M22").Value = -215
$ws.Range("N22").Value = -1375
$ws.Range("H31").Value = 3834.5264
$ws.Range("I31").Value = 1241
$ws.Range("J31").Value = 6048.512
$ws.Range("K31").Value = 1241
$ws.Range("L31").Value = 6048.512
$ws.Range("M31").Value = -946
$ws.Range("N31").Value = -6638.512
$ws.Range("H34").Value = 3834.5264
$ws.Range("I34").Value = 1241
$ws.Range("J34").Value = 6048.512
$ws.Range("K34").Value = 1241
$ws.Range("L34").Value = 6048.512
$ws.Range("M34").Value = -1039
$ws.Range("N34").Value = -6452.512
$ws.Range("H58").Value = 1160.5306
$ws.Range("I58").Value = 958.9666999999999
$ws.Range("J58").Value = 1478.7894
$ws.Range("K58").Value = 958.9666999999999
$ws.Range("L58").Value = 1478.7894
$ws.Range("M58").Value = -755.9666999999999
$ws.Range("N58").Value = -1884.7894
$ws.Range("H99").Value = 2041.3334
$ws.Range("I99").Value = 1763
$ws.Range("J99").Value = 2084.1538
$ws.Range("K99").Value = 1763
$ws.Range("L99").Value = 2084.1538
$ws.Range("M99").Value = -265
$ws.Range("N99").Value = -5080.1538
$ws.Range("H126").Value = 2041.3334
$ws.Range("I126").Value = 1763
$ws.Range("J126").Value = 2084.1538
$ws.Range("K126").Value = 5289
$ws.Range("L126").Value = 6252.4614
$ws.Range("M126").Value = -2819
$ws.Range("N126").Value = -11192.4614
$ws.Range("H132").Value = 3473687.5
$ws.Range("I132").Value = 1211.7222
$ws.Range("K132").Value = 3635.1666
$ws.Range("M132").Value = -1105.1666
$ws.Range("H134").Value = 3594.5833
$ws.Range("I134").Value = 3435.878
$ws.Range("J134").Value = 4524.143
$ws.Range("K134").Value = 10307.634
$ws.Range("L134").Value = 13572.429
$ws.Range("M134").Value = -7772.634
$ws.Range("N134").Value = -18642.429
$ws.Range("H136").Value = 1160.5306
$ws.Range("I136").Value = 958.9666999999999
$ws.Range("J136").Value = 1478.7894
$ws.Range("K136").Value = 2876.9001
$ws.Range("L136").Value = 4436.3682
$ws.Range("M136").Value = -326.9000999999998
$ws.Range("N136").Value = -9536.368200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 2924
$ws.Range("I93").Value = 823
$ws.Range("J93").Value = 3624.3333
$ws.Range("K93").Value = 2469
$ws.Range("L93").Value = 10872.9999
$ws.Range("M93").Value = -597
$ws.Range("N93").Value = -14616.9999
$ws.Range("H123").Value = 10000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 10000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 30000
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -34900
$ws.Range("H129").Value = 798584.4399999999
$ws.Range("I129").Value = 563.3333
$ws.Range("J129").Value = 1046246.2
$ws.Range("K129").Value = 1689.9999
$ws.Range("L129").Value = 3138738.6
$ws.Range("M129").Value = 3310.0001
$ws.Range("N129").Value = -3148738.6
$ws.Range("H130").Value = 1421.3334
$ws.Range("I130").Value = 686.25
$ws.Range("J130").Value = 2261.4285
$ws.Range("K130").Value = 2058.75
$ws.Range("L130").Value = 6784.2855
$ws.Range("M130").Value = 2961.25
$ws.Range("N130").Value = -16824.2855
$ws.Range("H131").Value = 3001.8245
$ws.Range("J131").Value = 3844.3257
$ws.Range("L131").Value = 11532.9771
$ws.Range("N131").Value = -21612.9771
$ws.Range("H133").Value = 13881.071
$ws.Range("I133").Value = 10606
$ws.Range("J133").Value = 15700.556
$ws.Range("K133").Value = 31818
$ws.Range("L133").Value = 47101.66800000001
$ws.Range("M133").Value = -26758
$ws.Range("N133").Value = -57221.66800000001
$ws.Range("H134").Value = 5148.7744
$ws.Range("I134").Value = 2376.6667
$ws.Range("J134").Value = 6899.579
$ws.Range("K134").Value = 7130.000100000001
$ws.Range("L134").Value = 20698.737
$ws.Range("M134").Value = -2060.000100000001
$ws.Range("N134").Value = -30838.737
$ws.Range("H136").Value = 4456.5
$ws.Range("I136").Value = 1843.3334
$ws.Range("J136").Value = 4917.647
$ws.Range("K136").Value = 5530.0002
$ws.Range("L136").Value = 14752.941
$ws.Range("M136").Value = -430.0002000000004
$ws.Range("N136").Value = -24952.941
$ws.Range("H137").Value = 6951312.5
$ws.Range("I137").Value = 16677889
$ws.Range("J137").Value = 3758.1428
$ws.Range("K137").Value = 50033667
$ws.Range("L137").Value = 11274.4284
$ws.Range("M137").Value = -50028567
$ws.Range("N137").Value = -21474.4284
$ws.Range("H138").Value = 7723.2
$ws.Range("I138").Value = 2260
$ws.Range("J138").Value = 13186.4
$ws.Range("K138").Value = 6780
$ws.Range("L138").Value = 39559.2
$ws.Range("M138").Value = -1640
$ws.Range("N138").Value = -49839.2
$ws.Range("H139").Value = 2433.5
$ws.Range("I139").Value = 1249.1666
$ws.Range("J139").Value = 3854.7
$ws.Range("K139").Value = 3747.4998
$ws.Range("L139").Value = 11564.1
$ws.Range("M139").Value = 1392.5002
$ws.Range("N139").Value = -21844.1
$ws.Range("H140").Value = 1634.4828
$ws.Range("I140").Value = 1352.381
$ws.Range("J140").Value = 2375
$ws.Range("K140").Value = 4057.143
$ws.Range("L140").Value = 7125
$ws.Range("M140").Value = 1122.857
$ws.Range("N140").Value = -17485
$ws.Range("H141").Value = 9092.182000000001
$ws.Range("I141").Value = 11154.375
$ws.Range("J141").Value = 7913.7856
$ws.Range("K141").Value = 33463.125
$ws.Range("L141").Value = 23741.3568
$ws.Range("M141").Value = -28283.125
$ws.Range("N141").Value = -34101.3568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2376.2097
$ws.Range("I132").Value = 1948.841
$ws.Range("K132").Value = 5846.522999999999
$ws.Range("M132").Value = -3316.522999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2084.1235
$ws.Range("I132").Value = 1933.0656
$ws.Range("J132").Value = 2544.85
$ws.Range("K132").Value = 5799.1968
$ws.Range("L132").Value = 7634.549999999999
$ws.Range("M132").Value = -3269.1968
$ws.Range("N132").Value = -12694.55
$ws.Range("H136").Value = 3877381
$ws.Range("I136").Value = 1179.6
$ws.Range("K136").Value = 3538.8
$ws.Range("M136").Value = -988.7999999999997

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4290650
$ws.Range("I132").Value = 1671.1389
$ws.Range("J132").Value = 9115752
$ws.Range("K132").Value = 5013.4167
$ws.Range("L132").Value = 27347256
$ws.Range("M132").Value = -2483.4167
$ws.Range("N132").Value = -27352316
$ws.Range("H136").Value = 1987.1632
$ws.Range("I136").Value = 1948.6323
$ws.Range("J136").Value = 2074.5
$ws.Range("K136").Value = 5845.8969
$ws.Range("L136").Value = 6223.5
$ws.Range("M136").Value = -3295.8969
$ws.Range("N136").Value = -11323.5
